$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.374922
$ws.Range("H2").Value = 4.124765999999999
$ws.Range("I2").Value = 0.2130738145062683
$ws.Range("J2").Value = 0.2130738145062683
$ws.Range("M2").Value = 5.209944
$ws.Range("N2").Value = 15.629832
$ws.Range("O2").Value = 0.1751928672265232
$ws.Range("P2").Value = 0.1751928672265232
$ws.Range("Q2").Value = 7.163266624367999
$ws.Range("R2").Value = 64.46939961931199
$ws.Range("S2").Value = 0.0373290124942455
$ws.Range("T2").Value = 0.03732901249424549
$ws.Range("G3").Value = 1.374922
$ws.Range("H3").Value = 4.124765999999999
$ws.Range("I3").Value = 0.2130738145062683
$ws.Range("J3").Value = 0.2130738145062683
$ws.Range("M3").Value = 6.497702
$ws.Range("O3").Value = 0.2184958310038485
$ws.Range("P3").Value = 0.2184958310038485
$ws.Range("Q3").Value = 8.933833429243998
$ws.Range("R3").Value = 80.40450086319599
$ws.Range("S3").Value = 0.04655574016570696
$ws.Range("T3").Value = 0.04655574016570695
$ws.Range("G4").Value = 1.374922
$ws.Range("H4").Value = 4.124765999999999
$ws.Range("I4").Value = 0.2130738145062683
$ws.Range("J4").Value = 0.2130738145062683
$ws.Range("M4").Value = 1.714656666666667
$ws.Range("N4").Value = 5.14397
$ws.Range("O4").Value = 0.057658127945791
$ws.Range("P4").Value = 0.057658127945791
$ws.Range("Q4").Value = 2.357519173446666
$ws.Range("R4").Value = 21.21767256102
$ws.Range("S4").Value = 0.01228543725870016
$ws.Range("T4").Value = 0.01228543725870016
$ws.Range("G5").Value = 1.374922
$ws.Range("H5").Value = 4.124765999999999
$ws.Range("I5").Value = 0.2130738145062683
$ws.Range("J5").Value = 0.2130738145062683
$ws.Range("M5").Value = 16.31603133333333
$ws.Range("N5").Value = 48.948094
$ws.Range("O5").Value = 0.5486531738238374
$ws.Range("P5").Value = 0.5486531738238374
$ws.Range("Q5").Value = 22.43327043288933
$ws.Range("R5").Value = 201.899433896004
$ws.Range("S5").Value = 0.1169036245876157
$ws.Range("T5").Value = 0.1169036245876157
$ws.Range("I6").Value = 0.4174838403400283
$ws.Range("J6").Value = 0.4174838403400283
$ws.Range("M6").Value = 5.209944
$ws.Range("N6").Value = 15.629832
$ws.Range("O6").Value = 0.1751928672265232
$ws.Range("P6").Value = 0.1751928672265232
$ws.Range("Q6").Value = 14.03526785612
$ws.Range("R6").Value = 126.31741070508
$ws.Range("S6").Value = 0.07314019100990959
$ws.Range("T6").Value = 0.07314019100990959
$ws.Range("I7").Value = 0.4174838403400283
$ws.Range("J7").Value = 0.4174838403400283
$ws.Range("M7").Value = 6.497702
$ws.Range("O7").Value = 0.2184958310038485
$ws.Range("P7").Value = 0.2184958310038485
$ws.Range("S7").Value = 0.09121847862577248
$ws.Range("T7").Value = 0.09121847862577248
$ws.Range("I8").Value = 0.4174838403400283
$ws.Range("J8").Value = 0.4174838403400283
$ws.Range("M8").Value = 1.714656666666667
$ws.Range("N8").Value = 5.14397
$ws.Range("O8").Value = 0.057658127945791
$ws.Range("P8").Value = 0.057658127945791
$ws.Range("Q8").Value = 4.61917932283889
$ws.Range("R8").Value = 41.57261390555001
$ws.Range("S8").Value = 0.02407133668162554
$ws.Range("T8").Value = 0.02407133668162554
$ws.Range("I9").Value = 0.4174838403400283
$ws.Range("J9").Value = 0.4174838403400283
$ws.Range("M9").Value = 16.31603133333333
$ws.Range("N9").Value = 48.948094
$ws.Range("O9").Value = 0.5486531738238374
$ws.Range("P9").Value = 0.5486531738238374
$ws.Range("Q9").Value = 43.95438225673444
$ws.Range("R9").Value = 395.58944031061
$ws.Range("S9").Value = 0.2290538340227207
$ws.Range("T9").Value = 0.2290538340227207
$ws.Range("G10").Value = 0.7672753333333334
$ws.Range("H10").Value = 2.301826
$ws.Range("I10").Value = 0.1189058594232269
$ws.Range("J10").Value = 0.1189058594232268
$ws.Range("M10").Value = 5.209944
$ws.Range("N10").Value = 15.629832
$ws.Range("O10").Value = 0.1751928672265232
$ws.Range("P10").Value = 0.1751928672265232
$ws.Range("Q10").Value = 3.997461519248001
$ws.Range("R10").Value = 35.977153673232
$ws.Range("S10").Value = 0.02083145844238901
$ws.Range("T10").Value = 0.02083145844238901
$ws.Range("G11").Value = 0.7672753333333334
$ws.Range("H11").Value = 2.301826
$ws.Range("I11").Value = 0.1189058594232269
$ws.Range("J11").Value = 0.1189058594232268
$ws.Range("M11").Value = 6.497702
$ws.Range("O11").Value = 0.2184958310038485
$ws.Range("P11").Value = 0.2184958310038485
$ws.Range("Q11").Value = 4.985526467950668
$ws.Range("R11").Value = 44.869738211556
$ws.Range("S11").Value = 0.02598043456590474
$ws.Range("T11").Value = 0.02598043456590473
$ws.Range("G12").Value = 0.7672753333333334
$ws.Range("H12").Value = 2.301826
$ws.Range("I12").Value = 0.1189058594232269
$ws.Range("J12").Value = 0.1189058594232268
$ws.Range("M12").Value = 1.714656666666667
$ws.Range("N12").Value = 5.14397
$ws.Range("O12").Value = 0.057658127945791
$ws.Range("P12").Value = 0.057658127945791
$ws.Range("Q12").Value = 1.315613765468889
$ws.Range("R12").Value = 11.84052388922
$ws.Range("S12").Value = 0.006855889256128653
$ws.Range("T12").Value = 0.006855889256128652
$ws.Range("G13").Value = 0.7672753333333334
$ws.Range("H13").Value = 2.301826
$ws.Range("I13").Value = 0.1189058594232269
$ws.Range("J13").Value = 0.1189058594232268
$ws.Range("M13").Value = 16.31603133333333
$ws.Range("N13").Value = 48.948094
$ws.Range("O13").Value = 0.5486531738238374
$ws.Range("P13").Value = 0.5486531738238374
$ws.Range("Q13").Value = 12.51888837996044
$ws.Range("R13").Value = 112.669995419644
$ws.Range("S13").Value = 0.06523807715880446
$ws.Range("T13").Value = 0.06523807715880445
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.616661
$ws.Range("H14").Value = 4.849983
$ws.Range("I14").Value = 0.2505364857304765
$ws.Range("J14").Value = 0.2505364857304765
$ws.Range("M14").Value = 5.209944
$ws.Range("N14").Value = 15.629832
$ws.Range("O14").Value = 0.1751928672265232
$ws.Range("P14").Value = 0.1751928672265232
$ws.Range("Q14").Value = 8.422713276984
$ws.Range("R14").Value = 75.804419492856
$ws.Range("S14").Value = 0.0438922052799791
$ws.Range("T14").Value = 0.0438922052799791
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.616661
$ws.Range("H15").Value = 4.849983
$ws.Range("I15").Value = 0.2505364857304765
$ws.Range("J15").Value = 0.2505364857304765
$ws.Range("M15").Value = 6.497702
$ws.Range("O15").Value = 0.2184958310038485
$ws.Range("P15").Value = 0.2184958310038485
$ws.Range("Q15").Value = 10.504581413022
$ws.Range("R15").Value = 94.541232717198
$ws.Range("S15").Value = 0.05474117764646429
$ws.Range("T15").Value = 0.05474117764646429
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.616661
$ws.Range("H16").Value = 4.849983
$ws.Range("I16").Value = 0.2505364857304765
$ws.Range("J16").Value = 0.2505364857304765
$ws.Range("M16").Value = 1.714656666666667
$ws.Range("N16").Value = 5.14397
$ws.Range("O16").Value = 0.057658127945791
$ws.Range("P16").Value = 0.057658127945791
$ws.Range("Q16").Value = 2.77201856139
$ws.Range("R16").Value = 24.94816705251
$ws.Range("S16").Value = 0.01444546474933666
$ws.Range("T16").Value = 0.01444546474933666
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.616661
$ws.Range("H17").Value = 4.849983
$ws.Range("I17").Value = 0.2505364857304765
$ws.Range("J17").Value = 0.2505364857304765
$ws.Range("M17").Value = 16.31603133333333
$ws.Range("N17").Value = 48.948094
$ws.Range("O17").Value = 0.5486531738238374
$ws.Range("P17").Value = 0.5486531738238374
$ws.Range("Q17").Value = 26.37749153137799
$ws.Range("R17").Value = 237.397423782402
$ws.Range("S17").Value = 0.1374576380546965
$ws.Range("T17").Value = 0.1374576380546965
